# "Hoàn thiện Ngoại Trú"
# Update the outpatient (Ngoại Trú) test-case reference record:
#  - "Data" sheet row 2: ticket/record number (A2) and ID card number (E2)
#  - "Check" sheet row 2: matching ticket/record number (A2)

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("A2").Value = 3012
$wsData.Range("E2").Value = 46200608012

$wsCheck = $wb.Worksheets.Item("Check")
$wsCheck.Range("A2").Value = 3012

# Restore focus to the "Data" sheet, matching the workbook's active sheet state.
$wsData.Activate()
$wsData.Range("Y5").Select()
